$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Top-20 hot stock rankings (A: 财联社, B: 东方财富, C: 同花顺)
$ws.Range("C2").Value = "闻泰科技"
$ws.Range("A3").Value = "香农芯创"
$ws.Range("B3").Value = "中国核建"
$ws.Range("C3").Value = "超颖电子"
$ws.Range("A4").Value = "N超颖"
$ws.Range("B4").Value = "香农芯创"
$ws.Range("C4").Value = "中信证券"
$ws.Range("B5").Value = "神开股份"
$ws.Range("C5").Value = "方正科技"
$ws.Range("A6").Value = "中际旭创"
$ws.Range("B6").Value = "中际旭创"
$ws.Range("C6").Value = "寒武纪"
$ws.Range("A7").Value = "中国卫星"
$ws.Range("B7").Value = "N超颖"
$ws.Range("C7").Value = "和而泰"
$ws.Range("B8").Value = "中国卫星"
$ws.Range("C8").Value = "香农芯创"
$ws.Range("A9").Value = "寒武纪-U"
$ws.Range("B9").Value = "三花智控"
$ws.Range("C9").Value = "大有能源"
$ws.Range("A10").Value = "科大国创"
$ws.Range("B10").Value = "平潭发展"
$ws.Range("C10").Value = "中际旭创"
$ws.Range("A11").Value = "方正科技"
$ws.Range("B11").Value = "科大国创"
$ws.Range("C11").Value = "珠江钢琴"
$ws.Range("A12").Value = "和而泰"
$ws.Range("B12").Value = "达华智能"
$ws.Range("C12").Value = "神开股份"
$ws.Range("A13").Value = "三花智控"
$ws.Range("B13").Value = "寒武纪-U"
$ws.Range("C13").Value = "中国核建"
$ws.Range("A14").Value = "平潭发展"
$ws.Range("B14").Value = "和而泰"
$ws.Range("C14").Value = "华建集团"
$ws.Range("A15").Value = "格尔软件"
$ws.Range("B15").Value = "神州信息"
$ws.Range("C15").Value = "上海电力"
$ws.Range("A16").Value = "闻泰科技"
$ws.Range("B16").Value = "方正科技"
$ws.Range("C16").Value = "三花智控"
$ws.Range("A17").Value = "达华智能"
$ws.Range("B17").Value = "山东墨龙"
$ws.Range("C17").Value = "大洋电机"
$ws.Range("A18").Value = "大洋电机"
$ws.Range("B18").Value = "立讯精密"
$ws.Range("C18").Value = "平潭发展"
$ws.Range("A19").Value = "航天科技"
$ws.Range("C19").Value = "国轩高科"
$ws.Range("A20").Value = "胜宏科技"
$ws.Range("B20").Value = "江波龙"
$ws.Range("C20").Value = "楚江新材"
$ws.Range("A21").Value = "中信证券"
$ws.Range("B21").Value = "大洋电机"
$ws.Range("C21").Value = "黄河旋风"

$wb.Save()
